$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 9).NumberFormat = "@"
$ws.Cells.Item(1, 9).Value = "11/03/2023"
$ws.Cells.Item(1, 9).NumberFormat = $ws.Cells.Item(1, 8).NumberFormat
$ws.Cells.Item(2, 3).Value = 1.147
$ws.Cells.Item(2, 4).Value = 1.177
$ws.Cells.Item(2, 5).Value = 17
$ws.Cells.Item(2, 6).Value = 9
$ws.Cells.Item(2, 7).Value = 13
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(2, 9).Value = 768.3
$ws.Cells.Item(2, 10).Value = -99.84680463360667
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 404
$ws.Cells.Item(3, 4).Value = 407
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 10
$ws.Cells.Item(3, 9).Value = 462
$ws.Cells.Item(3, 10).Value = -11.90476190476191
$ws.Cells.Item(4, 3).Value = 8
$ws.Cells.Item(4, 4).Value = 8
$ws.Cells.Item(4, 9).Value = 8
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(5, 3).Value = 300
$ws.Cells.Item(5, 4).Value = 312
$ws.Cells.Item(5, 5).Value = 11
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 9).Value = 257
$ws.Cells.Item(5, 10).Value = 21.40077821011672
$ws.Cells.Item(6, 3).Value = 67
$ws.Cells.Item(6, 4).Value = 67
$ws.Cells.Item(6, 9).Value = 56
$ws.Cells.Item(6, 10).Value = 19.64285714285714
$ws.Cells.Item(7, 3).Value = 54
$ws.Cells.Item(7, 4).Value = 55
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 7).Value = 1
$ws.Cells.Item(7, 9).Value = 109
$ws.Cells.Item(7, 10).Value = -49.54128440366973
$ws.Cells.Item(8, 3).Value = 294
$ws.Cells.Item(8, 4).Value = 295
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 7
$ws.Cells.Item(8, 9).Value = 17
$ws.Cells.Item(8, 10).Value = 1635.294117647059
$ws.Cells.Item(9, 3).Value = 53
$ws.Cells.Item(9, 4).Value = 53
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 4
$ws.Cells.Item(9, 9).Value = 40
$ws.Cells.Item(9, 10).Value = 32.49999999999999
$ws.Cells.Item(10, 2).Value = 3
$ws.Cells.Item(10, 3).Value = 200
$ws.Cells.Item(10, 4).Value = 248
$ws.Cells.Item(10, 5).Value = 49
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 9).Value = 537
$ws.Cells.Item(10, 10).Value = -53.81750465549349
$ws.Cells.Item(11, 3).Value = 333
$ws.Cells.Item(11, 4).Value = 338
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 3
$ws.Cells.Item(11, 9).Value = 239
$ws.Cells.Item(11, 10).Value = 41.42259414225941
$ws.Cells.Item(12, 2).Value = 5
$ws.Cells.Item(12, 3).Value = 285
$ws.Cells.Item(12, 4).Value = 643
$ws.Cells.Item(12, 5).Value = 211
$ws.Cells.Item(12, 6).Value = 5
$ws.Cells.Item(12, 7).Value = 5
$ws.Cells.Item(12, 8).Value = 137
$ws.Cells.Item(12, 9).Value = 676
$ws.Cells.Item(12, 10).Value = -4.881656804733725
$ws.Cells.Item(13, 3).Value = 20
$ws.Cells.Item(13, 4).Value = 29
$ws.Cells.Item(13, 5).Value = 7
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 9).Value = 722
$ws.Cells.Item(13, 10).Value = -95.98337950138504
$ws.Cells.Item(14, 2).Value = 4
$ws.Cells.Item(14, 3).Value = 354
$ws.Cells.Item(14, 4).Value = 622
$ws.Cells.Item(14, 5).Value = 195
$ws.Cells.Item(14, 7).Value = 3
$ws.Cells.Item(14, 8).Value = 59
$ws.Cells.Item(14, 9).Value = 825
$ws.Cells.Item(14, 10).Value = -24.60606060606061
$ws.Cells.Item(15, 3).Value = 236
$ws.Cells.Item(15, 4).Value = 247
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 9).Value = 256
$ws.Cells.Item(15, 10).Value = -3.515625
$ws.Cells.Item(16, 3).Value = 58
$ws.Cells.Item(16, 4).Value = 78
$ws.Cells.Item(16, 5).Value = 19
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 9).Value = 99
$ws.Cells.Item(16, 10).Value = -21.21212121212121
$ws.Cells.Item(17, 4).Value = 76
$ws.Cells.Item(17, 5).Value = 7
$ws.Cells.Item(17, 7).Value = 3
$ws.Cells.Item(17, 9).Value = 104
$ws.Cells.Item(17, 10).Value = -26.92307692307693
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 9).Value = 3
$ws.Cells.Item(18, 10).Value = -100
$ws.Cells.Item(19, 3).Value = 19
$ws.Cells.Item(19, 4).Value = 19
$ws.Cells.Item(19, 9).Value = 14
$ws.Cells.Item(19, 10).Value = 35.71428571428572
$ws.Cells.Item(20, 3).Value = 21
$ws.Cells.Item(20, 4).Value = 25
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 9).Value = 69
$ws.Cells.Item(20, 10).Value = -63.76811594202898
